# Refresh the Price (D) and Volume(1h) (E) columns of the cryptos
# table to the latest scraped snapshot. Most D-column prices are
# plain text that happens to look numeric (e.g. "228.87"), so a
# leading apostrophe forces Excel to keep them as text instead of
# auto-converting to a number, and the style is reset to Normal so
# no stray quote-prefix formatting is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '38.605.42'
$ws.Range('E2').Value = '  +2.46%  '
# Row 3
$ws.Range('D3').Value = '2.090.68'
$ws.Range('E3').Value = '  +2.91%  '
# Row 5
$ws.Range('D5').Value = '''228.87'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.55%  '
# Row 6
$ws.Range('D6').Value = '''0.613'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.22%  '
# Row 7
$ws.Range('D7').Value = '''61.26'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.18%  '
# Row 8
$ws.Range('E8').Value = '  +0.00%  '
# Row 9
$ws.Range('D9').Value = '''0.382'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.56%  '
# Row 10
$ws.Range('E10').Value = '  +2.65%  '
# Row 11
$ws.Range('D11').Value = '''0.105'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.38%  '
# Row 12
$ws.Range('D12').Value = '2.396.69'
$ws.Range('E12').Value = '  +2.79%  '
# Row 13
$ws.Range('D13').Value = '''14.80'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.26%  '
# Row 14
$ws.Range('D14').Value = '''22.27'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +6.19%  '
# Row 15
$ws.Range('E15').Value = '  +1.55%  '
# Row 16
$ws.Range('E16').Value = '  +5.06%  '
# Row 17
$ws.Range('D17').Value = '2.093.44'
$ws.Range('E17').Value = '  +3.18%  '
# Row 18
$ws.Range('D18').Value = '38.543.14'
$ws.Range('E18').Value = '  +2.40%  '
# Row 19
$ws.Range('D19').Value = '''71.05'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.23%  '
# Row 20
$ws.Range('E20').Value = '  +3.49%  '
# Row 21
$ws.Range('E21').Value = '  +1.47%  '
# Row 22
$ws.Range('D22').Value = '''226.17'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.19%  '
# Row 24
$ws.Range('D24').Value = '''2.40'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.47%  '
# Row 25
$ws.Range('E25').Value = '  +1.90%  '
# Row 26
$ws.Range('D26').Value = '''171.13'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.13%  '
# Row 27
$ws.Range('E27').Value = '  +0.94%  '
# Row 28
$ws.Range('E28').Value = '  +4.34%  '
# Row 29
$ws.Range('D29').Value = '''19.13'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.03%  '
# Row 30
$ws.Range('D30').Value = '''1.37'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +8.53%  '
# Row 31
$ws.Range('E31').Value = '  -0.01%  '
# Row 32
$ws.Range('D32').Value = '''2.32'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.94%  '
# Row 33
$ws.Range('D33').Value = '''4.79'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +6.54%  '
# Row 34
$ws.Range('D34').Value = '''4.48'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.63%  '
# Row 35
$ws.Range('E35').Value = '  +0.75%  '
# Row 36
$ws.Range('D36').Value = '''6.53'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.59%  '
# Row 37
$ws.Range('D37').Value = '''2.38'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.64%  '
# Row 38
$ws.Range('D38').Value = '''3.58'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +4.52%  '
# Row 39
$ws.Range('E39').Value = '  +0.05%  '
# Row 40
$ws.Range('D40').Value = '''18.58'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.48%  '
# Row 41
$ws.Range('D41').Value = '1.545.36'
$ws.Range('E41').Value = '  +0.47%  '
# Row 42
$ws.Range('D42').Value = '''99.82'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.01%  '
# Row 43
$ws.Range('D43').Value = '''0.0220'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.93%  '
# Row 44
$ws.Range('E44').Value = '  +1.14%  '
# Row 45
$ws.Range('D45').Value = '''0.0915'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.60%  '
# Row 46
$ws.Range('D46').Value = '''7.68'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +9.76%  '
# Row 47
$ws.Range('D47').Value = '''4.16'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.78%  '
# Row 48
$ws.Range('E48').Value = '  +1.37%  '
# Row 49
$ws.Range('E49').Value = '  +2.75%  '
# Row 50
$ws.Range('D50').Value = '''2.99'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.79%  '
# Row 51
$ws.Range('D51').Value = '2.287.66'
$ws.Range('E51').Value = '  +2.99%  '
